# Generate Report for Archive
#
# 1. Update status text from "Ready for handoff" to "In Translation"
#    wherever it appears (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# 2. Narrow the "Status" columns (Overview E:F, zh-cn C, de-de C) to
#    their new, tighter auto-fit width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Replace the status text -------------------------------------------------
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsZhCn.Range("C2:C4").Value     = "In Translation"
$wsDeDe.Range("C2:C4").Value     = "In Translation"

# --- Resize the affected columns ---------------------------------------------
# Target stored width is 13.4101845877511 "characters"; Excel's COM layer
# quantizes ColumnWidth to its internal pixel grid, so feed it the input
# (12.5) that lands on the nearest reachable grid value to the target.
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth     = 12.5
$wsDeDe.Range("C1").ColumnWidth     = 12.5
